# Generate Report for Handback
# Updates status of the 8abc4236... file (and its content-duplicate ab0d9c56...)
# from "Ready for handoff" to "Handed back: in sync with en-US" on all sheets,
# and records the new handback file / handback datetime on the zh-cn and de-de
# detail sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"
$ov.Range("E4").Value = "Handed back: in sync with en-US"
$ov.Range("F4").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("J3").Value = "8abc4236-22c1-4726-9f67-ea88bfe01212.13040e35b7035af807a3dcaaec0c5e4b1a6c31ad.zh-cn.xlf"
$zh.Range("K3").Value = "2016-10-14 07:59:16"

$zh.Range("C4").Value = "Handed back: in sync with en-US"
$zh.Range("J4").Value = "8abc4236-22c1-4726-9f67-ea88bfe01212.13040e35b7035af807a3dcaaec0c5e4b1a6c31ad.zh-cn.xlf"
$zh.Range("K4").Value = "2016-10-14 07:59:16"

# Rebuild hyperlinks in document order (A2, I2, A3, I3, A4, I4) so the new
# "Latest Handback File" links for rows 3 and 4 (column I) are inserted in
# the correct position.
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/16b10632a3620a9608c639457e416c1de16ef3f0/e2e/048ba33b-3a53-4b93-8c8f-5980f837820b.md", [Type]::Missing, [Type]::Missing, "048ba33b-3a53-4b93-8c8f-5980f837820b.md")
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/34d28387e6f5b1a11b6190ea1cc4e6d083f219a3/e2e/048ba33b-3a53-4b93-8c8f-5980f837820b.md", [Type]::Missing, [Type]::Missing, "048ba33b-3a53-4b93-8c8f-5980f837820b.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ed6d3cf5af8915831156ee801ab11318a396974/e2e/8abc4236-22c1-4726-9f67-ea88bfe01212.md", [Type]::Missing, [Type]::Missing, "8abc4236-22c1-4726-9f67-ea88bfe01212.md")
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5c1a6a4f7a2e9f9d4f2b6c7d8e9f0a1b2c3d4e5f/e2e/8abc4236-22c1-4726-9f67-ea88bfe01212.md", [Type]::Missing, [Type]::Missing, "8abc4236-22c1-4726-9f67-ea88bfe01212.md")
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ed6d3cf5af8915831156ee801ab11318a396974/e2e/ab0d9c56-bef2-49bf-8af1-5304e931b91d.md", [Type]::Missing, [Type]::Missing, "ab0d9c56-bef2-49bf-8af1-5304e931b91d.md")
$zh.Hyperlinks.Add($zh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5c1a6a4f7a2e9f9d4f2b6c7d8e9f0a1b2c3d4e5f/e2e/8abc4236-22c1-4726-9f67-ea88bfe01212.md", [Type]::Missing, [Type]::Missing, "8abc4236-22c1-4726-9f67-ea88bfe01212.md")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("J3").Value = "8abc4236-22c1-4726-9f67-ea88bfe01212.13040e35b7035af807a3dcaaec0c5e4b1a6c31ad.de-de.xlf"
$de.Range("K3").Value = "2016-10-14 07:59:33"

$de.Range("C4").Value = "Handed back: in sync with en-US"
$de.Range("J4").Value = "8abc4236-22c1-4726-9f67-ea88bfe01212.13040e35b7035af807a3dcaaec0c5e4b1a6c31ad.de-de.xlf"
$de.Range("K4").Value = "2016-10-14 07:59:33"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/16b10632a3620a9608c639457e416c1de16ef3f0/e2e/048ba33b-3a53-4b93-8c8f-5980f837820b.md", [Type]::Missing, [Type]::Missing, "048ba33b-3a53-4b93-8c8f-5980f837820b.md")
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/32dd3a75b636c9c93661b6437ffcc2f55c7b0220/e2e/048ba33b-3a53-4b93-8c8f-5980f837820b.md", [Type]::Missing, [Type]::Missing, "048ba33b-3a53-4b93-8c8f-5980f837820b.md")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ed6d3cf5af8915831156ee801ab11318a396974/e2e/8abc4236-22c1-4726-9f67-ea88bfe01212.md", [Type]::Missing, [Type]::Missing, "8abc4236-22c1-4726-9f67-ea88bfe01212.md")
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7d2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b/e2e/8abc4236-22c1-4726-9f67-ea88bfe01212.md", [Type]::Missing, [Type]::Missing, "8abc4236-22c1-4726-9f67-ea88bfe01212.md")
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ed6d3cf5af8915831156ee801ab11318a396974/e2e/ab0d9c56-bef2-49bf-8af1-5304e931b91d.md", [Type]::Missing, [Type]::Missing, "ab0d9c56-bef2-49bf-8af1-5304e931b91d.md")
$de.Hyperlinks.Add($de.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7d2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b/e2e/8abc4236-22c1-4726-9f67-ea88bfe01212.md", [Type]::Missing, [Type]::Missing, "8abc4236-22c1-4726-9f67-ea88bfe01212.md")
